$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Objetivos:" row (row 10) value with the new Portuguese objectives paragraph
$ws.Range("B10").Value = 'Este curso tem por objetivo fornecer aos alunos de Engenharia de Produção os princípios fundamentais da Química com enfoque tecnológico e nas aplicações industriais passíveis de serem encontradas na profissão.'
$ws.Range("C10").Value = 'Este curso tem por objetivo fornecer aos alunos de Engenharia de Produção os princípios fundamentais da Química com enfoque tecnológico e nas aplicações industriais passíveis de serem encontradas na profissão.'

# 2. Insert a new blank row at position 13 (everything from old row 13 down shifts to row 14+)
$ws.Rows("13").Insert()

# 3. New row 13 holds only the "Docentes responsaveis" value (no label in column A)
$ws.Range("A13").Clear()
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B13").Value = '198273 - Domingos Savio Giordani'
$ws.Range("C13").Value = '198273 - Domingos Savio Giordani'

# 4. Row 14 ("Programa resumido:") gets the new short summary text
$ws.Range("B14").Value = '1 – Conceitos básicos de Química; 2 – Os estados físicos da matéria e suas propriedades peculiares; 3 – Reações químicas; 4 – Noções de química orgânica; 5 – Materiais modernos; 6 – Tecnologia Química aplicada'
$ws.Range("C14").Value = '1 – Conceitos básicos de Química; 2 – Os estados físicos da matéria e suas propriedades peculiares; 3 – Reações químicas; 4 – Noções de química orgânica; 5 – Materiais modernos; 6 – Tecnologia Química aplicada'

# 5. Row 16 ("Programa:") gets the new full Portuguese syllabus text
$ws.Range("B16").Value = 'Programa em português1.Conceitos básicos de Química (2 horas)a.Estrutura Atômicab.Tabela Periódicac.Ligações Químicas2.Os estados físicos da matéria e suas propriedades peculiares (6 horas)a.O estado gasoso – pressão, relações PVT, gases ideais e reaisb.O estado líquido – soluções, forças intermoleculares, viscosidade, tensão superficial, pressão de vapor, mudanças de fasec.O estado sólido – classificação dos sólidos (moleculares, reticulares, metálicos e iônicos) 3.Reações químicas (8 horas)a.Tipos de reações (dupla-troca, oxirredução)b.Estequiometria em reações químicas (reagentes limitantes, pureza e rendimento)c.Energia e reações químicasd.Equilíbrio químico – soluções tampãoe.Fundamentos de corrosão4.Noções de química orgânica (6 horas)a.Hidrocarbonetos e suas principais propriedadesb.Combustíveis e combustãoc.Polímeros5.Tecnologia Química aplicada (8 horas)a.Papel e celuloseb.Açúcar e álcoolc.Sabões e detergentesd.Petróleo e gáse.Gases industriais f.Produção de vidros e cimento'
$ws.Range("C16").Value = 'Programa em português1.Conceitos básicos de Química (2 horas)a.Estrutura Atômicab.Tabela Periódicac.Ligações Químicas2.Os estados físicos da matéria e suas propriedades peculiares (6 horas)a.O estado gasoso – pressão, relações PVT, gases ideais e reaisb.O estado líquido – soluções, forças intermoleculares, viscosidade, tensão superficial, pressão de vapor, mudanças de fasec.O estado sólido – classificação dos sólidos (moleculares, reticulares, metálicos e iônicos) 3.Reações químicas (8 horas)a.Tipos de reações (dupla-troca, oxirredução)b.Estequiometria em reações químicas (reagentes limitantes, pureza e rendimento)c.Energia e reações químicasd.Equilíbrio químico – soluções tampãoe.Fundamentos de corrosão4.Noções de química orgânica (6 horas)a.Hidrocarbonetos e suas principais propriedadesb.Combustíveis e combustãoc.Polímeros5.Tecnologia Química aplicada (8 horas)a.Papel e celuloseb.Açúcar e álcoolc.Sabões e detergentesd.Petróleo e gáse.Gases industriais f.Produção de vidros e cimento'

# 6. Row 19 ("Metodo:") gets the evaluation method text
$ws.Range("B19").Value = 'Duas provas escritas e um seminário que, juntos, constituem a primeira avaliação.'
$ws.Range("C19").Value = 'Duas provas escritas e um seminário que, juntos, constituem a primeira avaliação.'

# 7. Row 20 ("Criterio:") gets the grading criteria text
$ws.Range("B20").Value = 'A nota de primeira avaliação será igual à média das notas das duas provas, com peso 7 somada à nota do seminário com peso 3. Alunos com nota de primeira avaliação igual ou superior a 5 estarão aprovados, com nota entre 3 e 4,9 em recuperação e abaixo de 3 reprovados.'
$ws.Range("C20").Value = 'A nota de primeira avaliação será igual à média das notas das duas provas, com peso 7 somada à nota do seminário com peso 3. Alunos com nota de primeira avaliação igual ou superior a 5 estarão aprovados, com nota entre 3 e 4,9 em recuperação e abaixo de 3 reprovados.'

# 8. Row 21 ("Norma de recuperacao:") gets the recovery norm text
$ws.Range("B21").Value = 'A recuperação se constituirá de uma prova abordando todos os assuntos do semestre, a nota de segunda avaliação será igual à média entre a nota de primeira avaliação e a prova de recuperação. Alunos com nota de segunda avaliação igual ou superior a 5 estarão aprovados e inferior a 5 reprovados.'
$ws.Range("C21").Value = 'A recuperação se constituirá de uma prova abordando todos os assuntos do semestre, a nota de segunda avaliação será igual à média entre a nota de primeira avaliação e a prova de recuperação. Alunos com nota de segunda avaliação igual ou superior a 5 estarão aprovados e inferior a 5 reprovados.'

# 9. New row 22 holds "Bibliografia:" label plus the bibliography text
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("B21").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("C21").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("A22").Value = 'Bibliografia:'
$ws.Range("B22").Value = 'BROWN, T.L. et al. Química a ciência central. 9.ed. São Paulo: Pearson Prentice Hall, 2005-2007ATKINS, P. Princípios de Química, questionando a vida moderna e o meio ambiente. 3ª Ed. Porto Alegre: Editora Bookman, 2006KOTZ, J. C. et al. Química geral e reações químicas, 9ª Edição, São Paulo, Cengage Learning, 2015.TOLENTINO, N. M. C. Processos Químicos Industriais, 1ª Edição, São Paulo, Érica, 2015.'
$ws.Range("C22").Value = 'BROWN, T.L. et al. Química a ciência central. 9.ed. São Paulo: Pearson Prentice Hall, 2005-2007ATKINS, P. Princípios de Química, questionando a vida moderna e o meio ambiente. 3ª Ed. Porto Alegre: Editora Bookman, 2006KOTZ, J. C. et al. Química geral e reações químicas, 9ª Edição, São Paulo, Cengage Learning, 2015.TOLENTINO, N. M. C. Processos Químicos Industriais, 1ª Edição, São Paulo, Érica, 2015.'
$ws.Rows("22").RowHeight = 120
